# "Add files via upload" — refresh of the Optical_Power claims sheet.
#
# The underlying CSV export was re-pulled: the six existing claim rows
# were re-ordered (their "OT" work-order numbers in column E got filled
# in, replacing the "Pendiente ADM" placeholder for most of them) and one
# brand-new claim (case 7092, SANTA FE AV. 2051) was appended as row 8.
#
# Net effect on the sheet: every data row (2-8) ends up holding a new set
# of values, and the used range grows from A1:N7 to A1:N8. The simplest,
# most robust way to reproduce that with Excel COM automation is to just
# (re)write every cell of rows 2-8 to its final value — equivalent to the
# end state, regardless of how the source list got reshuffled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (Caso), B (F. De Reclamo) and E (OT) hold values that look
# numeric/date-like ("7281", "9/22/2025", "01831884 ", ...) but must stay
# plain text (leading zeros / trailing spaces need to survive). Force
# text formatting on those columns for the data rows before writing so
# Excel doesn't silently convert them to numbers or date serials.
$ws.Range("A2:B8").NumberFormat = "@"
$ws.Range("E2:E8").NumberFormat = "@"

# --- Row 2: case 7281 (AZURDUY JUANA 2449) -------------------------------
$ws.Range("A2").Value = '7281'
$ws.Range("B2").Value = '9/22/2025'
$ws.Range("C2").Value = 'AZURDUY JUANA 2449'
$ws.Range("D2").Value = 13
$ws.Range("E2").Value = 'ICD30952422'
$ws.Range("F2").Value = 'Optical Power'
$ws.Range("G2").Value = 'Pendiente'
$ws.Range("H2").Value = 'Cable en panza y cortado'
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = '{"direccionesNormalizadas": [{"altura": 2449, "cod_calle": 1151, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.467279", "y": "-34.551117"}, "direccion": "AZURDUY JUANA 2449, CABA", "nombre_calle": "AZURDUY JUANA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K2").Value = -58.467279
$ws.Range("L2").Value = -34.551117
$ws.Range("M2").Value = 'Saavedra'
$ws.Range("N2").Value = 'Capital Norte'

# --- Row 3: case 4757 (GARAY, JUAN DE AV. 819) ---------------------------
$ws.Range("A3").Value = '4757 '
$ws.Range("B3").Value = '12/11/2025'
$ws.Range("C3").Value = 'GARAY, JUAN DE AV. 819'
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = '01831884 '
$ws.Range("F3").Value = 'Optical Power'
$ws.Range("G3").Value = 'Pendiente'
$ws.Range("H3").Value = 'tendido bajo'
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = '{"direccionesNormalizadas": [{"altura": 819, "cod_calle": 7026, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.376986", "y": "-34.625210"}, "direccion": "GARAY, JUAN DE AV. 819, CABA", "nombre_calle": "GARAY, JUAN DE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K3").Value = -58.376986
$ws.Range("L3").Value = -34.62521
$ws.Range("M3").Value = 'San Telmo'
$ws.Range("N3").Value = 'Capital Sur'

# --- Row 4: case 4756 (GARAY, JUAN DE AV. 799) ---------------------------
$ws.Range("A4").Value = '4756 '
$ws.Range("B4").Value = '12/11/2025'
$ws.Range("C4").Value = 'GARAY, JUAN DE AV. 799'
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = '01831840 '
$ws.Range("F4").Value = 'Optical Power'
$ws.Range("G4").Value = 'Pendiente'
$ws.Range("H4").Value = 'tendido bajo'
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = '{"direccionesNormalizadas": [{"altura": 799, "cod_calle": 7026, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.376455", "y": "-34.624886"}, "direccion": "GARAY, JUAN DE AV. 799, CABA", "nombre_calle": "GARAY, JUAN DE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K4").Value = -58.376455
$ws.Range("L4").Value = -34.624886
$ws.Range("M4").Value = 'San Telmo'
$ws.Range("N4").Value = 'Capital Sur'

# --- Row 5: case S00299847 (VALLE, ARISTOBULO DEL 1357) ------------------
$ws.Range("A5").Value = 'S00299847'
$ws.Range("B5").Value = '12/15/2025'
$ws.Range("C5").Value = 'VALLE, ARISTOBULO DEL 1357'
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = '01939674 '
$ws.Range("F5").Value = 'Optical Power'
$ws.Range("G5").Value = 'Pendiente'
$ws.Range("H5").Value = 'tendido bajo'
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = '{"direccionesNormalizadas": [{"altura": 1357, "cod_calle": 4057, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.371495", "y": "-34.635725"}, "direccion": "VALLE, ARISTOBULO DEL 1357, CABA", "nombre_calle": "VALLE, ARISTOBULO DEL", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K5").Value = -58.371495
$ws.Range("L5").Value = -34.635725
$ws.Range("M5").Value = 'San Telmo'
$ws.Range("N5").Value = 'Capital Sur'

# --- Row 6: case 7937 (BOULOGNE MER 323) ---------------------------------
$ws.Range("A6").Value = '7937'
$ws.Range("B6").Value = '12/19/2025'
$ws.Range("C6").Value = 'BOULOGNE MER 323'
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = '01939828 '
$ws.Range("F6").Value = 'Optical Power'
$ws.Range("G6").Value = 'Pendiente'
$ws.Range("H6").Value = 'cables colgando'
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = '{"direccionesNormalizadas": [{"altura": 323, "cod_calle": 2106, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.407174", "y": "-34.606292"}, "direccion": "BOULOGNE SUR MER 323, CABA", "nombre_calle": "BOULOGNE SUR MER", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K6").Value = -58.407174
$ws.Range("L6").Value = -34.606292
$ws.Range("M6").Value = 'Almagro'
$ws.Range("N6").Value = 'Capital Sur'

# --- Row 7: case 5004 (SANTA FE AV. 4830) --------------------------------
$ws.Range("A7").Value = '5004'
$ws.Range("B7").Value = '12/19/2025'
$ws.Range("C7").Value = 'SANTA FE AV. 4830'
$ws.Range("D7").Value = 14
$ws.Range("E7").Value = '01939765 '
$ws.Range("F7").Value = 'Optical Power'
$ws.Range("G7").Value = 'Pendiente'
$ws.Range("H7").Value = 'cables sueltos cortados'
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = '{"direccionesNormalizadas": [{"altura": 4830, "cod_calle": 20057, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.427883", "y": "-34.577829"}, "direccion": "SANTA FE AV. 4830, CABA", "nombre_calle": "SANTA FE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K7").Value = -58.427883
$ws.Range("L7").Value = -34.577829
$ws.Range("M7").Value = 'Palermo'
$ws.Range("N7").Value = 'Capital Sur'

# --- Row 8 (NEW): case 7092 (SANTA FE AV. 2051) --------------------------
$ws.Range("A8").Value = '7092'
$ws.Range("B8").Value = '12/26/2025'
$ws.Range("C8").Value = 'SANTA FE AV. 2051'
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 'Pendiente ADM'
$ws.Range("F8").Value = 'Optical Power'
$ws.Range("G8").Value = 'Pendiente'
$ws.Range("H8").Value = 'tendido bajo'
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = '{"direccionesNormalizadas": [{"altura": 2051, "cod_calle": 20057, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.396730", "y": "-34.595570"}, "direccion": "SANTA FE AV. 2051, CABA", "nombre_calle": "SANTA FE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K8").Value = -58.39673
$ws.Range("L8").Value = -34.59557
$ws.Range("M8").Value = 'Recoleta'
$ws.Range("N8").Value = 'Capital Sur'
